$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gdrive Links")

# Fix the typo'd Google Drive id in the shared string used by C2 (and,
# transitively, the cached formula result in D2).
$ws.Range("C2").Value = "https://drive.google.com/open?id=1Tik7O5yXSrohqb0jOe80EC4tXYmzJKRi"

# Turn C2 into a real hyperlink pointing at the same Google Drive URL.
$ws.Hyperlinks.Add($ws.Range("C2"), "https://drive.google.com/open?id=1Tik7O5yXSrohqb0jOe80EC4tXYmzJKRi")

# Re-apply the workbook's existing "Hyperlink" cell style so C2 matches
# the style already used elsewhere in the sheet instead of a freshly
# minted duplicate.
$ws.Range("C2").Style = "Hyperlink"

# Move the active selection from D7 to D2.
$ws.Activate()
$ws.Range("D2").Select() | Out-Null
